$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the register in A10 from DVC_SAMPLE_DATA to DVC_FLUSH_SAMPLE_DATA
$ws.Range("A10").Value = "DVC_FLUSH_SAMPLE_DATA"

# Update the selection to reflect the user having selected the whole row 10
$ws.Range("A10:XFD10").Select()
